# Auto-generated Excel COM-interop script
# Updates leve-crafting price/profit columns (H-N) per scheduled market-data refresh
$wb = $excel.ActiveWorkbook

# Sheet ALC, row 109 (hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Sheet ALC, row 116 (hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 8141925.5
$ws.Range("I116").Value = 9227316
$ws.Range("K116").Value = 9227316
$ws.Range("M116").Value = -9223874

# Sheet ALC, row 137 (hunk 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 43480148
$ws.Range("I137").Value = 66667760
$ws.Range("J137").Value = 3373.25
$ws.Range("K137").Value = 200003280
$ws.Range("L137").Value = 10119.75
$ws.Range("M137").Value = -200000730
$ws.Range("N137").Value = -15219.75

# Sheet ARM, row 32 (hunk 3)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4409.061
$ws.Range("I32").Value = 2404.5789
$ws.Range("K32").Value = 2404.5789
$ws.Range("M32").Value = -2117.5789

# Sheet ARM, row 61 (hunk 4)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2859.4243
$ws.Range("I61").Value = 1711.5454
$ws.Range("J61").Value = 5155.1816
$ws.Range("K61").Value = 1711.5454
$ws.Range("L61").Value = 5155.1816
$ws.Range("M61").Value = -1499.5454
$ws.Range("N61").Value = -5579.1816

# Sheet ARM, row 74 (hunk 5)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5450.967
$ws.Range("I74").Value = 1584.762
$ws.Range("J74").Value = 14472.111
$ws.Range("K74").Value = 1584.762
$ws.Range("L74").Value = 14472.111
$ws.Range("M74").Value = -710.7619999999999
$ws.Range("N74").Value = -16220.111

# Sheet ARM, row 77 (hunk 6)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5450.967
$ws.Range("I77").Value = 1584.762
$ws.Range("J77").Value = 14472.111
$ws.Range("K77").Value = 7923.809999999999
$ws.Range("L77").Value = 72360.55500000001
$ws.Range("M77").Value = -3555.809999999999
$ws.Range("N77").Value = -81096.55500000001

# Sheet ARM, row 132 (hunk 7)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3666.5
$ws.Range("I132").Value = 3103.5715
$ws.Range("J132").Value = 4229.4287
$ws.Range("K132").Value = 9310.7145
$ws.Range("L132").Value = 12688.2861
$ws.Range("M132").Value = -6780.7145
$ws.Range("N132").Value = -17748.2861

# Sheet ARM, row 133 (hunk 8)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 36166.668
$ws.Range("J133").Value = 36166.668
$ws.Range("L133").Value = 36166.668
$ws.Range("N133").Value = -41226.668

# Sheet ARM, row 136 (hunk 9)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2859.4243
$ws.Range("I136").Value = 1711.5454
$ws.Range("J136").Value = 5155.1816
$ws.Range("K136").Value = 5134.6362
$ws.Range("L136").Value = 15465.5448
$ws.Range("M136").Value = -2584.6362
$ws.Range("N136").Value = -20565.5448

# Sheet ARM, row 139 (hunk 10)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 60702.332
$ws.Range("J139").Value = 60702.332
$ws.Range("L139").Value = 60702.332
$ws.Range("N139").Value = -70982.33199999999

# Sheet BSM, row 59 (hunk 11)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 45000
$ws.Range("J59").Value = 45000
$ws.Range("L59").Value = 45000
$ws.Range("N59").Value = -46694

# Sheet BSM, row 81 (hunk 12)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 39319.668
$ws.Range("J81").Value = 39319.668
$ws.Range("L81").Value = 39319.668
$ws.Range("N81").Value = -41441.668

# Sheet BSM, row 84 (hunk 13)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 39319.668
$ws.Range("J84").Value = 39319.668
$ws.Range("L84").Value = 117959.004
$ws.Range("N84").Value = -128567.004

# Sheet BSM, row 107 (hunk 14)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 705
$ws.Range("I107").Value = 697.1429000000001
$ws.Range("J107").Value = 714.1667
$ws.Range("K107").Value = 697.1429000000001
$ws.Range("L107").Value = 714.1667
$ws.Range("M107").Value = 1222.8571
$ws.Range("N107").Value = -4554.1667

# Sheet BSM, row 134 (hunk 15)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2676.1777
$ws.Range("I134").Value = 1682.7
$ws.Range("J134").Value = 4663.1333
$ws.Range("K134").Value = 5048.1
$ws.Range("L134").Value = 13989.3999
$ws.Range("M134").Value = -2513.1
$ws.Range("N134").Value = -19059.3999

# Sheet CRP, row 31 (hunk 16)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1357
$ws.Range("I31").Value = 1104.6364
$ws.Range("J31").Value = 2745
$ws.Range("K31").Value = 1104.6364
$ws.Range("L31").Value = 2745
$ws.Range("M31").Value = -809.6364000000001
$ws.Range("N31").Value = -3335

# Sheet CRP, row 34 (hunk 17)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1357
$ws.Range("I34").Value = 1104.6364
$ws.Range("J34").Value = 2745
$ws.Range("K34").Value = 1104.6364
$ws.Range("L34").Value = 2745
$ws.Range("M34").Value = -902.6364000000001
$ws.Range("N34").Value = -3149

# Sheet CRP, row 58 (hunk 18)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3129
$ws.Range("I58").Value = 2193.6
$ws.Range("J58").Value = 3752.6
$ws.Range("K58").Value = 2193.6
$ws.Range("L58").Value = 3752.6
$ws.Range("M58").Value = -1990.6
$ws.Range("N58").Value = -4158.6

# Sheet CRP, row 122 (hunk 19)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2338.9167
$ws.Range("I122").Value = 1307.3334
$ws.Range("J122").Value = 5433.6665
$ws.Range("K122").Value = 3922.0002
$ws.Range("L122").Value = 16300.9995
$ws.Range("M122").Value = -1472.0002
$ws.Range("N122").Value = -21200.9995

# Sheet CRP, row 132 (hunk 20)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3430.2917
$ws.Range("I132").Value = 2641.5833
$ws.Range("J132").Value = 4219
$ws.Range("K132").Value = 7924.749899999999
$ws.Range("L132").Value = 12657
$ws.Range("M132").Value = -5394.749899999999
$ws.Range("N132").Value = -17717

# Sheet CRP, row 134 (hunk 21)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3918.9375
$ws.Range("I134").Value = 1891.6666
$ws.Range("J134").Value = 5135.3
$ws.Range("K134").Value = 5674.9998
$ws.Range("L134").Value = 15405.9
$ws.Range("M134").Value = -3139.9998
$ws.Range("N134").Value = -20475.9

# Sheet CRP, row 136 (hunk 22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3129
$ws.Range("I136").Value = 2193.6
$ws.Range("J136").Value = 3752.6
$ws.Range("K136").Value = 6580.799999999999
$ws.Range("L136").Value = 11257.8
$ws.Range("M136").Value = -4030.799999999999
$ws.Range("N136").Value = -16357.8

# Sheet GSM, row 102 (hunk 23)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1428.25
$ws.Range("I102").Value = 1208.2222
$ws.Range("J102").Value = 1711.1428
$ws.Range("K102").Value = 1208.2222
$ws.Range("L102").Value = 1711.1428
$ws.Range("M102").Value = 413.7778000000001
$ws.Range("N102").Value = -4955.1428

# Sheet GSM, row 132 (hunk 24)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2839.6226
$ws.Range("I132").Value = 2421.8386
$ws.Range("J132").Value = 3428.318
$ws.Range("K132").Value = 7265.5158
$ws.Range("L132").Value = 10284.954
$ws.Range("M132").Value = -4735.5158
$ws.Range("N132").Value = -15344.954

# Sheet GSM, row 137 (hunk 25)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# Sheet GSM, row 138 (hunk 26)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Sheet LTW, row 46 (hunk 27)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 790.4167
$ws.Range("I46").Value = 698.3
$ws.Range("J46").Value = 1251
$ws.Range("K46").Value = 698.3
$ws.Range("L46").Value = 1251
$ws.Range("M46").Value = -510.3
$ws.Range("N46").Value = -1627

# Sheet LTW, row 61 (hunk 28)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5572.577
$ws.Range("I61").Value = 5117.4707
$ws.Range("J61").Value = 6432.222
$ws.Range("K61").Value = 5117.4707
$ws.Range("L61").Value = 6432.222
$ws.Range("M61").Value = -4915.4707
$ws.Range("N61").Value = -6836.222

# Sheet LTW, row 113 (hunk 29)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5572.577
$ws.Range("I113").Value = 5117.4707
$ws.Range("J113").Value = 6432.222
$ws.Range("K113").Value = 5117.4707
$ws.Range("L113").Value = 6432.222
$ws.Range("M113").Value = -2947.4707
$ws.Range("N113").Value = -10772.222

# Sheet LTW, row 132 (hunk 30)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3517.3022
$ws.Range("I132").Value = 2624.2693
$ws.Range("J132").Value = 4883.1177
$ws.Range("K132").Value = 7872.8079
$ws.Range("L132").Value = 14649.3531
$ws.Range("M132").Value = -5342.8079
$ws.Range("N132").Value = -19709.3531

# Sheet LTW, row 136 (hunk 31)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4745.759
$ws.Range("I136").Value = 2338.0557
$ws.Range("J136").Value = 8685.637000000001
$ws.Range("K136").Value = 7014.1671
$ws.Range("L136").Value = 26056.911
$ws.Range("M136").Value = -4464.1671
$ws.Range("N136").Value = -31156.911

# Sheet WVR, row 122 (hunk 32)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 54442.473
$ws.Range("I122").Value = 101770.7
$ws.Range("J122").Value = 1855.5555
$ws.Range("K122").Value = 305312.1
$ws.Range("L122").Value = 5566.666499999999
$ws.Range("M122").Value = -302862.1
$ws.Range("N122").Value = -10466.6665

# Sheet WVR, row 132 (hunk 33)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 20002874
$ws.Range("I132").Value = 26318132
$ws.Range("K132").Value = 78954396
$ws.Range("M132").Value = -78951866

# Sheet WVR, row 136 (hunk 34)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9037703
$ws.Range("I136").Value = 15921514
$ws.Range("J136").Value = 2700.375
$ws.Range("K136").Value = 47764542
$ws.Range("L136").Value = 8101.125
$ws.Range("M136").Value = -47761992
$ws.Range("N136").Value = -13201.125
